$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Platform Coverage")

# Insert a new row above row 2, pushing the existing age-band rows down by one.
$ws.Rows("2:2").Insert()

# Populate the newly inserted row 2 with a new MDA age-band entry (ages 5-15).
$ws.Range("A2").Value = "All"
$ws.Range("B2").Value = "Treatment"
$ws.Range("C2").Value = "Campaign"
$ws.Range("D2").Value = "MDA"
$ws.Range("F2").Value = 5
$ws.Range("G2").Value = 15
$ws.Range("H2").Value = 0.6
$ws.Range("J2").Value = 0.6
$ws.Range("L2").Value = 0.6
$ws.Range("N2").Value = 0.6

# The row-insert shifted row 2's old H/J/L/N coverage values down into row 3
# along with everything else; row 3 (ages 2-15) should only carry the P:AD
# coverage figures, so clear out the stray H/J/L/N cells it inherited.
$ws.Range("H3").ClearContents()
$ws.Range("J3").ClearContents()
$ws.Range("L3").ClearContents()
$ws.Range("N3").ClearContents()
